# Generate Report for Handoff
# Updates the localization-status report to reflect that b.md has been
# handed off again (new handoff file generated). Its status is no longer
# "Handed back: in sync with en-US" but "Ready for handoff", since the
# handback version is stale relative to the newly published source.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6f7b54689d416143c8f3c689dd07a6a6ce0ca06/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e37f560dd2ea9ed2de1679cb89d02493e10e4795/e2e/b.md."

# ---- Overview sheet: row 3 corresponds to b.md ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-26 12:38:01"

# ---- zh-cn sheet: row 3 corresponds to b.md ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe keeps this a text "False" (matching the other text
# True/False cells in the sheet) instead of Excel auto-typing it as a
# native boolean; resetting the style afterwards drops the quote-prefix
# formatting flag so the cell style matches the rest of the column.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-26 12:37:56"
$wsZhCn.Range("P3").Value = $errorDetail
# 39.17 character-units renders as the stored OOXML column width of 40
# (same ratio already used by the other width=40 columns in this sheet).
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet: row 3 corresponds to b.md ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-26 12:38:01"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
